$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, pushing existing row 34 (and below) down by one.
$ws.Rows.Item(34).Insert()

# Populate the new row 34 with the new data entry (same market/category metadata,
# new date and new volume/price figures).
$ws.Cells.Item(34, 1).Value = 5
$ws.Cells.Item(34, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(34, 3).Value = "Maule"
$ws.Cells.Item(34, 4).Value = 44414
$ws.Cells.Item(34, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(34, 5).Value = 7
$ws.Cells.Item(34, 6).Value = 100112009
$ws.Cells.Item(34, 7).Value = "Acelga"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 400
$ws.Cells.Item(34, 11).Value = 2000
$ws.Cells.Item(34, 12).Value = 2000
$ws.Cells.Item(34, 13).Value = 2000
$ws.Cells.Item(34, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(34, 15).Value = "Región del Maule"
$ws.Cells.Item(34, 16).Value = 500
$ws.Cells.Item(34, 17).Value = 4
$ws.Cells.Item(34, 18).Value = "Hortaliza"
